$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Replace "mike" with "nameofuserprofile" inside the path, keeping formatting.
$d.Content.Find.Execute("c:\users\mike\documents\scripts\myFunction.ps1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "c:\users\nameofuserprofile\documents\scripts\myFunction.ps1", 2)
